$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.01642412646811
$ws.Range("C2").Value = 0.397697122200384
$ws.Range("E2").Value = 0.1434344796069061
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 0.3482483691223237
$ws.Range("H2").Value = 0.5312191985089427
$ws.Range("L2").Value = 0.187685287646012
$ws.Range("M2").Value = 0.2130373181354841
$ws.Range("O2").Value = 1.695424113027528

$ws.Range("B3").Value = 0.9052719095669772
$ws.Range("C3").Value = 0.386533692453412
$ws.Range("E3").Value = 0.145521296762097
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 0.355185999970459
$ws.Range("H3").Value = 0.5388386736436885
$ws.Range("L3").Value = 0.1849634506432523
$ws.Range("M3").Value = 0.1957259501154809
$ws.Range("O3").Value = 1.725824444783484

$ws.Range("B4").Value = 0.8368595318659118
$ws.Range("C4").Value = 0.3797393369040662
$ws.Range("E4").Value = 0.146878514165989
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 0.3598434649888915
$ws.Range("H4").Value = 0.5438455276547955
$ws.Range("L4").Value = 0.1833815078344117
$ws.Range("M4").Value = 0.1851061718096787
$ws.Range("O4").Value = 1.74601247791513

$ws.Range("B5").Value = 0.8089415427449467
$ws.Range("C5").Value = 0.3769860179966713
$ws.Range("E5").Value = 0.1474506911843021
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.3618411684525782
$ws.Range("H5").Value = 0.5459684202143364
$ws.Range("L5").Value = 0.1827593691587879
$ws.Range("M5").Value = 0.1807812096036869
$ws.Range("O5").Value = 1.754621369473909

$ws.Range("B6").Value = 0.8043034610475956
$ws.Range("C6").Value = 0.376529773846471
$ws.Range("E6").Value = 0.1475468547240014
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.3621789023932038
$ws.Range("H6").Value = 0.546325910190042
$ws.Range("L6").Value = 0.1826574253747779
$ws.Range("M6").Value = 0.1800632239389728
$ws.Range("O6").Value = 1.756073928629462

$ws.Range("B7").Value = 0.8364831769713987
$ws.Range("C7").Value = 0.379702141729922
$ws.Range("E7").Value = 0.1468861533925199
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.3598700031891866
$ws.Range("H7").Value = 0.5438738234974636
$ws.Range("L7").Value = 0.1833730262000728
$ws.Range("M7").Value = 0.1850478325889
$ws.Range("O7").Value = 1.746127034166008

$ws.Range("B8").Value = 0.9781341190138164
$ws.Range("C8").Value = 0.3938357102232715
$ws.Range("E8").Value = 0.1441382687121811
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.3505577570554479
$ws.Range("H8").Value = 0.5337781993219934
$ws.Range("L8").Value = 0.1867283076238664
$ws.Range("M8").Value = 0.2070665677236079
$ws.Range("O8").Value = 1.705589851522404

$ws.Range("B9").Value = 1.254533533221604
$ws.Range("C9").Value = 0.4220146491887533
$ws.Range("E9").Value = 0.1393512690111021
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.3354640258774353
$ws.Range("H9").Value = 0.516588400464812
$ws.Range("L9").Value = 0.1940139906361935
$ws.Range("M9").Value = 0.2503093371135137
$ws.Range("O9").Value = 1.638201985529378

$ws.Range("B10").Value = 1.456686851178517
$ws.Range("C10").Value = 0.4429845164912933
$ws.Range("E10").Value = 0.1362000049274494
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.3263216607411508
$ws.Range("H10").Value = 0.5055505426671587
$ws.Range("L10").Value = 0.1997947818802288
$ws.Range("M10").Value = 0.2821072673195886
$ws.Range("O10").Value = 1.59610984010186

$ws.Range("B11").Value = 1.548437453016049
$ws.Range("C11").Value = 0.4525790983999798
$ws.Range("E11").Value = 0.1348456150258757
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.3225887326768273
$ws.Range("H11").Value = 0.5008750856532913
$ws.Range("L11").Value = 0.2025170967008307
$ws.Range("M11").Value = 0.2965766990427667
$ws.Range("O11").Value = 1.578579784243288

$ws.Range("B12").Value = 1.583149148501036
$ws.Range("C12").Value = 0.4562199653991286
$ws.Range("E12").Value = 0.1343441085881603
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.3212367060141474
$ws.Range("H12").Value = 0.4991543670064544
$ws.Range("L12").Value = 0.2035612280134984
$ws.Range("M12").Value = 0.3020562721202964
$ws.Range("O12").Value = 1.572174943303992

$ws.Range("B13").Value = 1.575674831717834
$ws.Range("C13").Value = 0.4554355066397591
$ws.Range("E13").Value = 0.1344516113438445
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.3215251468909983
$ws.Range("H13").Value = 0.4995227398037585
$ws.Range("L13").Value = 0.2033357674247043
$ws.Range("M13").Value = 0.3008761391301391
$ws.Range("O13").Value = 1.573543947295818

$ws.Range("B14").Value = 1.551293863584874
$ws.Range("C14").Value = 0.4528784842133291
$ws.Range("E14").Value = 0.1348041279308708
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.3224762653386932
$ws.Range("H14").Value = 0.5007325231701785
$ws.Range("L14").Value = 0.2026027328478364
$ws.Range("M14").Value = 0.2970275023553768
$ws.Range("O14").Value = 1.578048171873007

$ws.Range("B15").Value = 1.536355553303395
$ws.Range("C15").Value = 0.4513132146979331
$ws.Range("E15").Value = 0.1350215352031832
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.323066877538615
$ws.Range("H15").Value = 0.5014800344637962
$ws.Range("L15").Value = 0.2021554515416994
$ws.Range("M15").Value = 0.2946701325894381
$ws.Range("O15").Value = 1.580837559012636

$ws.Range("B16").Value = 1.450686330994074
$ws.Range("C16").Value = 0.4423585710762268
$ws.Range("E16").Value = 0.1362901092224329
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.3265742135856655
$ws.Range("H16").Value = 0.5058630567966134
$ws.Range("L16").Value = 0.1996187299396581
$ws.Range("M16").Value = 0.2811617170658423
$ws.Range("O16").Value = 1.597288091507608

$ws.Range("B17").Value = 1.398075781790055
$ws.Range("C17").Value = 0.4368791018263778
$ws.Range("E17").Value = 0.1370886011735578
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.3288351660608342
$ws.Range("H17").Value = 0.5086404899182853
$ws.Range("L17").Value = 0.1980862040852998
$ws.Range("M17").Value = 0.2728756397297829
$ws.Range("O17").Value = 1.607794880380652

$ws.Range("B18").Value = 1.367795946284161
$ws.Range("C18").Value = 0.4337326798971901
$ws.Range("E18").Value = 0.1375553217494909
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.3301756837246685
$ws.Range("H18").Value = 0.5102705331250732
$ws.Range("L18").Value = 0.1972134560360672
$ws.Range("M18").Value = 0.2681101396712933
$ws.Range("O18").Value = 1.613990343928322

$ws.Range("B19").Value = 1.357540404707379
$ws.Range("C19").Value = 0.4326682631985648
$ws.Range("E19").Value = 0.1377146249571702
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.330636434264413
$ws.Range("H19").Value = 0.5108280241587906
$ws.Range("L19").Value = 0.19691945828977
$ws.Range("M19").Value = 0.2664967072507238
$ws.Range("O19").Value = 1.616114144884094

$ws.Range("B20").Value = 1.403678312217778
$ws.Range("C20").Value = 0.4374618625268738
$ws.Range("E20").Value = 0.1370028294441837
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.3285903337723539
$ws.Range("H20").Value = 0.5083414595038818
$ws.Range("L20").Value = 0.1982484420443171
$ws.Range("M20").Value = 0.2737576643630888
$ws.Range("O20").Value = 1.606660653323402

$ws.Range("B21").Value = 1.558456037079338
$ws.Range("C21").Value = 0.4536293399414149
$ws.Range("E21").Value = 0.1347002767052189
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.3221952259305425
$ws.Range("H21").Value = 0.5003758290902951
$ws.Range("L21").Value = 0.2028176838071971
$ws.Range("M21").Value = 0.2981579342931369
$ws.Range("O21").Value = 1.576718831173622

$ws.Range("B22").Value = 1.659423204060033
$ws.Range("C22").Value = 0.4642399043168837
$ws.Range("E22").Value = 0.1332617042318747
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.3183745628815089
$ws.Range("H22").Value = 0.4954599868755594
$ws.Range("L22").Value = 0.2058811464469272
$ws.Range("M22").Value = 0.314106601129474
$ws.Range("O22").Value = 1.55851096678262

$ws.Range("B23").Value = 1.605553103914474
$ws.Range("C23").Value = 0.4585729159505547
$ws.Range("E23").Value = 0.1340234355511842
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.3203807848205642
$ws.Range("H23").Value = 0.4980570947196696
$ws.Range("L23").Value = 0.2042390770263154
$ws.Range("M23").Value = 0.3055944484641486
$ws.Range("O23").Value = 1.568104078243053

$ws.Range("B24").Value = 1.401145511873949
$ws.Range("C24").Value = 0.437198384559224
$ws.Range("E24").Value = 0.1370415829918168
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.3287008958010631
$ws.Range("H24").Value = 0.5084765475561071
$ws.Range("L24").Value = 0.1981750683381733
$ws.Range("M24").Value = 0.2733589063720956
$ws.Range("O24").Value = 1.60717295471116

$ws.Range("B25").Value = 1.179915728125422
$ws.Range("C25").Value = 0.4143435418645822
$ws.Range("E25").Value = 0.1405819784845279
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.3392065140422247
$ws.Range("H25").Value = 0.5209593418374894
$ws.Range("L25").Value = 0.1919676871614584
$ws.Range("M25").Value = 0.2386053687456027
$ws.Range("O25").Value = 1.655132284126296
